{"js": "// Replace the three-digit x one-digit multiplication equations throughout\n// the document body (table cells) with the newly generated set of\n// equations, preserving run formatting (font, size, etc.) by doing an\n// in-place text replace on each matched range.\nconst replacements = [\n  [\"291\u00d73=873\", \"992\u00d74=3968\"],\n  [\"730\u00d78=5840\", \"383\u00d74=1532\"],\n  [\"931\u00d78=7448\", \"174\u00d73=522\"],\n  [\"137\u00d77=959\", \"169\u00d74=676\"],\n  [\"982\u00d73=2946\", \"887\u00d78=7096\"],\n  [\"493\u00d75=2465\", \"459\u00d72=918\"],\n  [\"452\u00d75=2260\", \"677\u00d78=5416\"],\n  [\"665\u00d78=5320\", \"272\u00d79=2448\"],\n  [\"289\u00d72=578\", \"789\u00d77=5523\"],\n  [\"461\u00d72=922\", \"943\u00d72=1886\"],\n  [\"472\u00d78=3776\", \"835\u00d79=7515\"],\n  [\"995\u00d76=5970\", \"264\u00d75=1320\"],\n  [\"894\u00d75=4470\", \"401\u00d74=1604\"],\n  [\"866\u00d73=2598\", \"809\u00d76=4854\"],\n  [\"299\u00d72=598\", \"536\u00d75=2680\"],\n  [\"178\u00d78=1424\", \"933\u00d79=8397\"],\n  [\"945\u00d73=2835\", \"473\u00d75=2365\"],\n  [\"421\u00d77=2947\", \"136\u00d79=1224\"],\n  [\"106\u00d74=424\", \"207\u00d75=1035\"],\n  [\"266\u00d72=532\", \"837\u00d78=6696\"],\n  [\"675\u00d76=4050\", \"716\u00d79=6444\"],\n  [\"499\u00d78=3992\", \"706\u00d77=4942\"],\n  [\"814\u00d73=2442\", \"573\u00d76=3438\"],\n  [\"673\u00d79=6057\", \"269\u00d77=1883\"],\n  [\"490\u00d73=1470\", \"942\u00d72=1884\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication equations throughout\n# the document (table cells) with the newly generated set of equations,\n# using Find/Replace on the document's Content range so that run\n# formatting (font, size, etc.) is preserved for the surrounding text.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"291\u00d73=873\", \"992\u00d74=3968\"),\n    @(\"730\u00d78=5840\", \"383\u00d74=1532\"),\n    @(\"931\u00d78=7448\", \"174\u00d73=522\"),\n    @(\"137\u00d77=959\", \"169\u00d74=676\"),\n    @(\"982\u00d73=2946\", \"887\u00d78=7096\"),\n    @(\"493\u00d75=2465\", \"459\u00d72=918\"),\n    @(\"452\u00d75=2260\", \"677\u00d78=5416\"),\n    @(\"665\u00d78=5320\", \"272\u00d79=2448\"),\n    @(\"289\u00d72=578\", \"789\u00d77=5523\"),\n    @(\"461\u00d72=922\", \"943\u00d72=1886\"),\n    @(\"472\u00d78=3776\", \"835\u00d79=7515\"),\n    @(\"995\u00d76=5970\", \"264\u00d75=1320\"),\n    @(\"894\u00d75=4470\", \"401\u00d74=1604\"),\n    @(\"866\u00d73=2598\", \"809\u00d76=4854\"),\n    @(\"299\u00d72=598\", \"536\u00d75=2680\"),\n    @(\"178\u00d78=1424\", \"933\u00d79=8397\"),\n    @(\"945\u00d73=2835\", \"473\u00d75=2365\"),\n    @(\"421\u00d77=2947\", \"136\u00d79=1224\"),\n    @(\"106\u00d74=424\", \"207\u00d75=1035\"),\n    @(\"266\u00d72=532\", \"837\u00d78=6696\"),\n    @(\"675\u00d76=4050\", \"716\u00d79=6444\"),\n    @(\"499\u00d78=3992\", \"706\u00d77=4942\"),\n    @(\"814\u00d73=2442\", \"573\u00d76=3438\"),\n    @(\"673\u00d79=6057\", \"269\u00d77=1883\"),\n    @(\"490\u00d73=1470\", \"942\u00d72=1884\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
